$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")
$ws.Activate()

# 1. Update the text value of C2 (shared string change)
$ws.Range("C2").Value = "Cook,Bedroom,Bathroom,2,1,2,Every week  Every day  No repeat  Every month"

# 2. Match C2's style/formatting to B2's (so the old unused style can be pruned)
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Set column C width
$ws.Columns.Item(3).ColumnWidth = 41.04

# 4. Update the active selection to C3
$ws.Range("C3").Select() | Out-Null
